$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Update the cached date text in the "Datum:" TIME field
#    7. November 2016  ->  21. November 2016
# ------------------------------------------------------------------
$d.Content.Find.Execute("7. November 2016", $false, $false, $false, $false, $false, $true, 1, $false, "21. November 2016", 2) | Out-Null

# ------------------------------------------------------------------
# 2. Extend the "Fetcher" heading to "Fetcher (Gery)" and relocate
#    the "_GoBack" bookmark from the end of the
#    "... einander folgen, da die Interessen gleichartig sind."
#    paragraph into the middle of the new heading text (between
#    "(Gery" and ")"), mirroring where the cursor was left after
#    typing the addition.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$rng = $d.Content
$rng.Find.Execute("Fetcher") | Out-Null
$rng.Collapse(0)
$headingEnd = $rng.Start
$rng.InsertAfter(" (Gery)")

# Range covering " (Gery" only, to compute the boundary right after it.
$split = $d.Range($headingEnd, $headingEnd)
$split.MoveEnd(1, 6) | Out-Null

$bmRng = $d.Range($split.End, $split.End)
$d.Bookmarks.Add("_GoBack", $bmRng) | Out-Null
